# Update odds values for the Jogos da Semana FlashScore workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 updates
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.88

# Row 6 updates
$ws.Range("H6").Value = 4.3
$ws.Range("I6").Value = 6.5
$ws.Range("J6").Value = 1.83
$ws.Range("K6").Value = 2.37
$ws.Range("L6").Value = 6.1
$ws.Range("P6").Value = 4.15
$ws.Range("S6").Value = 1.29
$ws.Range("T6").Value = 3.32
$ws.Range("W6").Value = 6.8
$ws.Range("Y6").Value = 7
$ws.Range("Z6").Value = 8
$ws.Range("AD6").Value = 7.6
$ws.Range("AK6").Value = 110
$ws.Range("AO6").Value = 6.2
$ws.Range("AQ6").Value = 17
$ws.Range("AT6").Value = 3.2
$ws.Range("AU6").Value = 7.8
